$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(54, 8).Value = 4720
$ws.Cells.Item(54, 9).Value = 2980
$ws.Cells.Item(54, 11).Value = 2980
$ws.Cells.Item(54, 13).Value = -2494
$ws.Cells.Item(137, 8).Value = 843.6316
$ws.Cells.Item(137, 9).Value = 658.0625
$ws.Cells.Item(137, 11).Value = 1974.1875
$ws.Cells.Item(137, 13).Value = 575.8125
$ws.Cells.Item(139, 8).Value = 59220
$ws.Cells.Item(139, 10).Value = 59220
$ws.Cells.Item(139, 12).Value = 59220
$ws.Cells.Item(139, 14).Value = -69500
$ws.Cells.Item(140, 8).Value = 83232.73
$ws.Cells.Item(140, 10).Value = 83232.73
$ws.Cells.Item(140, 12).Value = 83232.73
$ws.Cells.Item(140, 14).Value = -93592.73

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(37, 8).Value = 6960.6
$ws.Cells.Item(37, 10).Value = 8200.75
$ws.Cells.Item(37, 12).Value = 8200.75
$ws.Cells.Item(37, 14).Value = -8746.75
$ws.Cells.Item(61, 8).Value = 3100.0967
$ws.Cells.Item(61, 9).Value = 3507.0908
$ws.Cells.Item(61, 11).Value = 3507.0908
$ws.Cells.Item(61, 13).Value = -3295.0908
$ws.Cells.Item(74, 8).Value = 1656.6136
$ws.Cells.Item(74, 9).Value = 1565.6765
$ws.Cells.Item(74, 10).Value = 1965.8
$ws.Cells.Item(74, 11).Value = 1565.6765
$ws.Cells.Item(74, 12).Value = 1965.8
$ws.Cells.Item(74, 13).Value = -691.6765
$ws.Cells.Item(74, 14).Value = -3713.8
$ws.Cells.Item(77, 8).Value = 1656.6136
$ws.Cells.Item(77, 9).Value = 1565.6765
$ws.Cells.Item(77, 10).Value = 1965.8
$ws.Cells.Item(77, 11).Value = 7828.3825
$ws.Cells.Item(77, 12).Value = 9829
$ws.Cells.Item(77, 13).Value = -3460.3825
$ws.Cells.Item(77, 14).Value = -18565
$ws.Cells.Item(122, 8).Value = 3837.6296
$ws.Cells.Item(122, 9).Value = 3979.7827
$ws.Cells.Item(122, 10).Value = 3020.25
$ws.Cells.Item(122, 11).Value = 11939.3481
$ws.Cells.Item(122, 12).Value = 9060.75
$ws.Cells.Item(122, 13).Value = -9489.348100000001
$ws.Cells.Item(122, 14).Value = -13960.75
$ws.Cells.Item(136, 8).Value = 3100.0967
$ws.Cells.Item(136, 9).Value = 3507.0908
$ws.Cells.Item(136, 11).Value = 10521.2724
$ws.Cells.Item(136, 13).Value = -7971.2724
$ws.Cells.Item(139, 8).Value = 47500
$ws.Cells.Item(139, 10).Value = 47500
$ws.Cells.Item(139, 12).Value = 47500
$ws.Cells.Item(139, 14).Value = -57780

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 5321.974
$ws.Cells.Item(134, 9).Value = 1130.08
$ws.Cells.Item(134, 10).Value = 12807.5
$ws.Cells.Item(134, 11).Value = 3390.24
$ws.Cells.Item(134, 12).Value = 38422.5
$ws.Cells.Item(134, 13).Value = -855.2399999999998
$ws.Cells.Item(134, 14).Value = -43492.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3369.6206
$ws.Cells.Item(31, 9).Value = 2282.8262
$ws.Cells.Item(31, 10).Value = 7535.6665
$ws.Cells.Item(31, 11).Value = 2282.8262
$ws.Cells.Item(31, 12).Value = 7535.6665
$ws.Cells.Item(31, 13).Value = -1987.8262
$ws.Cells.Item(31, 14).Value = -8125.6665
$ws.Cells.Item(34, 8).Value = 3369.6206
$ws.Cells.Item(34, 9).Value = 2282.8262
$ws.Cells.Item(34, 10).Value = 7535.6665
$ws.Cells.Item(34, 11).Value = 2282.8262
$ws.Cells.Item(34, 12).Value = 7535.6665
$ws.Cells.Item(34, 13).Value = -2080.8262
$ws.Cells.Item(34, 14).Value = -7939.6665
$ws.Cells.Item(50, 8).Value = 8429.143
$ws.Cells.Item(50, 10).Value = 8429.143
$ws.Cells.Item(50, 12).Value = 8429.143
$ws.Cells.Item(50, 14).Value = -9679.143
$ws.Cells.Item(51, 8).Value = 9400.75
$ws.Cells.Item(51, 10).Value = 9400.75
$ws.Cells.Item(51, 12).Value = 9400.75
$ws.Cells.Item(51, 14).Value = -10872.75
$ws.Cells.Item(60, 8).Value = 8300.429
$ws.Cells.Item(60, 10).Value = 8300.429
$ws.Cells.Item(60, 12).Value = 8300.429
$ws.Cells.Item(60, 14).Value = -9322.429
$ws.Cells.Item(61, 8).Value = 9400.75
$ws.Cells.Item(61, 10).Value = 9400.75
$ws.Cells.Item(61, 12).Value = 9400.75
$ws.Cells.Item(61, 14).Value = -10096.75
$ws.Cells.Item(68, 8).Value = 17700.5
$ws.Cells.Item(68, 10).Value = 17700.5
$ws.Cells.Item(68, 12).Value = 17700.5
$ws.Cells.Item(68, 14).Value = -19198.5
$ws.Cells.Item(71, 8).Value = 17700.5
$ws.Cells.Item(71, 10).Value = 17700.5
$ws.Cells.Item(71, 12).Value = 53101.5
$ws.Cells.Item(71, 14).Value = -60589.5
$ws.Cells.Item(74, 8).Value = 16863.7
$ws.Cells.Item(74, 10).Value = 16863.7
$ws.Cells.Item(74, 12).Value = 16863.7
$ws.Cells.Item(74, 14).Value = -18611.7
$ws.Cells.Item(77, 8).Value = 16863.7
$ws.Cells.Item(77, 10).Value = 16863.7
$ws.Cells.Item(77, 12).Value = 50591.10000000001
$ws.Cells.Item(77, 14).Value = -59327.10000000001
$ws.Cells.Item(122, 8).Value = 1714.9354
$ws.Cells.Item(122, 9).Value = 1353.5238
$ws.Cells.Item(122, 10).Value = 2473.9
$ws.Cells.Item(122, 11).Value = 4060.5714
$ws.Cells.Item(122, 12).Value = 7421.700000000001
$ws.Cells.Item(122, 13).Value = -1610.5714
$ws.Cells.Item(122, 14).Value = -12321.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 822.76
$ws.Cells.Item(131, 10).Value = 826.2857
$ws.Cells.Item(131, 12).Value = 2478.8571
$ws.Cells.Item(131, 14).Value = -12558.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1887.6428
$ws.Cells.Item(122, 9).Value = 1795.6364
$ws.Cells.Item(122, 10).Value = 2225
$ws.Cells.Item(122, 11).Value = 5386.9092
$ws.Cells.Item(122, 12).Value = 6675
$ws.Cells.Item(122, 13).Value = -2936.9092
$ws.Cells.Item(122, 14).Value = -11575
$ws.Cells.Item(132, 8).Value = 1675.7693
$ws.Cells.Item(132, 9).Value = 1210.125
$ws.Cells.Item(132, 10).Value = 2420.8
$ws.Cells.Item(132, 11).Value = 3630.375
$ws.Cells.Item(132, 12).Value = 7262.400000000001
$ws.Cells.Item(132, 13).Value = -1100.375
$ws.Cells.Item(132, 14).Value = -12322.4
$ws.Cells.Item(138, 8).Value = 67999.336
$ws.Cells.Item(138, 10).Value = 67999.336
$ws.Cells.Item(138, 12).Value = 67999.336
$ws.Cells.Item(138, 14).Value = -78279.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 958.3043
$ws.Cells.Item(55, 9).Value = 1419.9
$ws.Cells.Item(55, 10).Value = 603.2308
$ws.Cells.Item(55, 11).Value = 1419.9
$ws.Cells.Item(55, 12).Value = 603.2308
$ws.Cells.Item(55, 13).Value = -1246.9
$ws.Cells.Item(55, 14).Value = -949.2308
$ws.Cells.Item(136, 8).Value = 3074.4473
$ws.Cells.Item(136, 9).Value = 1910.3214
$ws.Cells.Item(136, 10).Value = 6334
$ws.Cells.Item(136, 11).Value = 5730.9642
$ws.Cells.Item(136, 12).Value = 19002
$ws.Cells.Item(136, 13).Value = -3180.9642
$ws.Cells.Item(136, 14).Value = -24102
$ws.Cells.Item(139, 8).Value = 79733.336
$ws.Cells.Item(139, 10).Value = 79733.336
$ws.Cells.Item(139, 12).Value = 79733.336
$ws.Cells.Item(139, 14).Value = -90013.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1545.9231
$ws.Cells.Item(122, 9).Value = 1464.2727
$ws.Cells.Item(122, 11).Value = 4392.8181
$ws.Cells.Item(122, 13).Value = -1942.8181
$ws.Cells.Item(136, 8).Value = 904.7907
$ws.Cells.Item(136, 9).Value = 660.93335
$ws.Cells.Item(136, 10).Value = 1467.5385
$ws.Cells.Item(136, 11).Value = 1982.80005
$ws.Cells.Item(136, 12).Value = 4402.6155
$ws.Cells.Item(136, 13).Value = 567.1999499999999
$ws.Cells.Item(136, 14).Value = -9502.6155
$ws.Cells.Item(141, 8).Value = 67233.336
